# Generate Report for Handoff
# Adds a new "handback" row for 64d947ee-a767-493c-8a90-ea0403e5866b to each
# sheet of the localization-status workbook (Overview, zh-cn, de-de),
# mirroring the existing rows for 571c7103-8660-4f11-8e8b-df8803d0e27d.

$wb = $excel.ActiveWorkbook

$newId     = "64d947ee-a767-493c-8a90-ea0403e5866b"
$newFile   = "$newId.md"
$newToken  = "bd75e25f15b28351a1a39d6513bbb031ca67d6f8"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/$newFile",
    "",
    "",
    $newFile
)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-21 22:35:35"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/$newFile",
    "",
    "",
    $newFile
)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/86c5f1c3fcf4eeb680ca7d9a65b3ee542c89a9bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newId.$newToken.zh-cn.xlf",
    "",
    "",
    "$newId.$newToken.zh-cn.xlf"
)
$wsZh.Range("E3").Value = "2016-03-21 22:35:31"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("J3").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/243d6bdb350f89f1f294ac662ae3930b5ff8cbfa/e2e/$newFile",
    "",
    "",
    $newFile
)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34ae5c474cf2361ac996412ee2a82e4e64ab8941/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newId.$newToken.de-de.xlf",
    "",
    "",
    "$newId.$newToken.de-de.xlf"
)
$wsDe.Range("E3").Value = "2016-03-21 22:35:35"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("J3").Value = "Include"
